$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Glossary")

# Add the new glossary entries that were added for the search/get-term
# recognition feature (rows 28-32 below the existing terms).
$ws.Range("A28").Value = "getting"
$ws.Range("B28").Value = "erste"
$ws.Range("C28").Value = "Not approved"

$ws.Range("A29").Value = "finding"
$ws.Range("B29").Value = "Geeigneten"
$ws.Range("C29").Value = "Approved"

$ws.Range("A30").Value = "register"
$ws.Range("B30").Value = "registrieren|anmelden"
$ws.Range("C30").Value = "Approved|Approved"

$ws.Range("A31").Value = "obtain"
$ws.Range("B31").Value = "erhalten"
$ws.Range("C31").Value = "Approved"

$ws.Range("A32").Value = "discount"
$ws.Range("B32").Value = "Rabatt|Diskont"
$ws.Range("C32").Value = "Approved|Approved"

$ws.Range("D32").Select() | Out-Null
